$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the page count for "Researching Information Systems and Computing" (row 11, column C)
$ws.Range("C11").Value = 178

# Update selection to B15
$ws.Range("B15").Select()
